# EPBDS-4097 Deployment isolation feature, fix bug in class loader isolation
#
# The "beans" data source used by project2's TestBean demo is renamed to
# "beans2" (to isolate it from project1's same-named data source), and a
# second "javabean" example method (printJavaBeanSecond) is added to the
# sheet, mirroring the existing printJavaBean layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "beans" data source to "beans2" -----------------------
$ws.Range("F4").Value = "Data TestBean beans2"
$ws.Range("B12").Value = 'return "project2"+print(beans2[0]);'

# --- Insert two new rows (18,19) before the old "Environment" block,
#     which pushes it from rows 19-21 down to rows 21-23 ---------------
$ws.Rows("18:19").Insert()

# Row 18: "Method String printJavaBeanSecond()" label, styled/merged the
# same way as the other method-name rows (e.g. B11:D11, B16:D16)
$ws.Range("B18").Value = "Method String printJavaBeanSecond()"
$ws.Range("B18:D18").Merge() | Out-Null
foreach ($addr in "B18","C18","D18") {
    $cell = $ws.Range($addr)
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
    $cell.HorizontalAlignment = -4108
}

# Row 19: matching body/formula row
$ws.Range("B19").Value = 'return "project2"+print(beans2[0]);'
$ws.Range("B19:D19").Merge() | Out-Null
foreach ($addr in "B19","C19","D19") {
    $cell = $ws.Range($addr)
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
    $cell.HorizontalAlignment = -4108
}

# --- The "import" example at the bottom (now rows 21-23) now points at
#     org.openl.example2 -------------------------------------------------
$ws.Range("C23").Value = "org.openl.example2"

# --- Selection moves onto the renamed data-source cell -----------------
$ws.Range("F4:G4").Select() | Out-Null
